$wb = $excel.ActiveWorkbook

# The file 9e95c7cc-1260-4d05-953e-b03aa7e01f94.md has moved from
# "Ready for handoff" to "In Translation" for both locales, so update the
# per-locale status cells on the Overview sheet and the Status column on
# each locale's report sheet.

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E4").Value = "In Translation"
$overview.Range("F4").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C4").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C4").Value = "In Translation"
